$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay text (avoid Excel auto-numeric coercion)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '28.380.73'
$ws.Range("E2").Value = '  -0.61%  '

# Row 3
$ws.Range("D3").Value = '1.833.22'
$ws.Range("E3").Value = '  +2.36%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.42%  '

# Row 5
$ws.Range("D5").Value = '330.13'
$ws.Range("E5").Value = '  +0.31%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.26%  '

# Row 7
$ws.Range("D7").Value = '0.4469'
$ws.Range("E7").Value = '  +1.80%  '

# Row 8
$ws.Range("D8").Value = '0.3789'
$ws.Range("E8").Value = '  +1.39%  '

# Row 9
$ws.Range("D9").Value = '44.85'
$ws.Range("E9").Value = '  -1.31%  '

# Row 10
$ws.Range("D10").Value = '0.07784'
$ws.Range("E10").Value = '  +2.58%  '

# Row 11
$ws.Range("D11").Value = '1.141'
$ws.Range("E11").Value = '  +0.96%  '

# Row 12
$ws.Range("D12").Value = '22.41'
$ws.Range("E12").Value = '  -0.80%  '

# Row 13
$ws.Range("D13").Value = '1.000'
$ws.Range("E13").Value = '  -0.43%  '

# Row 14
$ws.Range("D14").Value = '6.382'
$ws.Range("E14").Value = '  +2.75%  '

# Row 15
$ws.Range("D15").Value = '7.576'
$ws.Range("E15").Value = '  +1.30%  '

# Row 16
$ws.Range("D16").Value = '1.837.99'
$ws.Range("E16").Value = '  +2.53%  '

# Row 17
$ws.Range("D17").Value = '93.38'
$ws.Range("E17").Value = '  +16.16%  '

# Row 18
$ws.Range("D18").Value = '0.00001088'
$ws.Range("E18").Value = '  +0.07%  '

# Row 19
$ws.Range("D19").Value = '0.06430'
$ws.Range("E19").Value = '  -3.98%  '

# Row 20
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.24%  '

# Row 21
$ws.Range("D21").Value = '17.64'
$ws.Range("E21").Value = '  +0.90%  '

# Row 22
$ws.Range("D22").Value = '6.380'
$ws.Range("E22").Value = '  +2.91%  '

# Row 23
$ws.Range("D23").Value = '0.5411'
$ws.Range("E23").Value = '  +0.64%  '

# Row 24
$ws.Range("D24").Value = '28.410.50'
$ws.Range("E24").Value = '  -0.49%  '

# Row 25
$ws.Range("D25").Value = '11.79'
$ws.Range("E25").Value = '  +1.00%  '

# Row 26
$ws.Range("D26").Value = '2.242'
$ws.Range("E26").Value = '  -7.96%  '

# Row 27
$ws.Range("D27").Value = '20.91'
$ws.Range("E27").Value = '  +2.75%  '

# Row 28
$ws.Range("D28").Value = '154.84'
$ws.Range("E28").Value = '  +1.38%  '

# Row 29
$ws.Range("D29").Value = '2.389'
$ws.Range("E29").Value = '  +2.54%  '

# Row 30
$ws.Range("D30").Value = '2.042.61'
$ws.Range("E30").Value = '  +2.17%  '

# Row 31
$ws.Range("D31").Value = '129.06'
$ws.Range("E31").Value = '  -1.08%  '

# Row 32
$ws.Range("D32").Value = '1.217'
$ws.Range("E32").Value = '  -6.38%  '

# Row 33
$ws.Range("D33").Value = '5.962'
$ws.Range("E33").Value = '  +3.15%  '

# Row 34
$ws.Range("D34").Value = '0.09312'
$ws.Range("E34").Value = '  +0.80%  '

# Row 35
$ws.Range("D35").Value = '3.673'
$ws.Range("E35").Value = '  -7.64%  '

# Row 36
$ws.Range("D36").Value = '13.26'
$ws.Range("E36").Value = '  +9.90%  '

# Row 37
$ws.Range("D37").Value = '0.02365'
$ws.Range("E37").Value = '  +2.16%  '

# Row 38
$ws.Range("D38").Value = '0.2204'
$ws.Range("E38").Value = '  -1.60%  '

# Row 39
$ws.Range("D39").Value = '5.231'
$ws.Range("E39").Value = '  +0.83%  '

# Row 40
$ws.Range("D40").Value = '0.6661'
$ws.Range("E40").Value = '  +1.39%  '

# Row 41
$ws.Range("D41").Value = '0.06280'
$ws.Range("E41").Value = '  +0.48%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.211'
$ws.Range("E42").Value = '  +2.81%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.202'
$ws.Range("E43").Value = '  +0.44%  '

# Row 44
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.28%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '1.404'
$ws.Range("E45").Value = '  -1.31%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '14.04'
$ws.Range("E46").Value = '  +1.14%  '

# Row 47
$ws.Range("D47").Value = '0.6169'
$ws.Range("E47").Value = '  +1.72%  '

# Row 48
$ws.Range("D48").Value = '3.790'
$ws.Range("E48").Value = '  -0.57%  '

# Row 49
$ws.Range("D49").Value = '2.073'
$ws.Range("E49").Value = '  +3.30%  '

# Row 50
$ws.Range("D50").Value = '127.96'
$ws.Range("E50").Value = '  +0.49%  '

# Row 51
$ws.Range("D51").Value = '0.07019'
$ws.Range("E51").Value = '  +0.32%  '
